$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty row 8, shifting rows 9-14 up to become rows 8-13
$ws.Rows.Item(8).Delete()

# Update the active selection to match the post-edit state
$ws.Range("B16").Select()
